$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.926.63"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "1.902.49"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'0.8002"
$ws.Range("E5").Value = "  +5.85%  "
$ws.Range("D6").Value = "'240.84"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Value = "'0.9999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'0.3123"
$ws.Range("E8").Value = "  +2.64%  "
$ws.Range("D9").Value = "'26.27"
$ws.Range("E9").Value = "  +3.13%  "
$ws.Range("D10").Value = "'0.07067"
$ws.Range("E10").Value = "  +3.42%  "
$ws.Range("D11").Value = "'0.07976"
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").Value = "1.912.28"
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("D13").Value = "'0.7376"
$ws.Range("E13").Value = "  -1.03%  "
$ws.Range("D14").Value = "'5.179"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").Value = "'92.53"
$ws.Range("E15").Value = "  +1.58%  "
$ws.Range("D16").Value = "29.914.00"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "'13.93"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "'5.870"
$ws.Range("E18").Value = "  -1.11%  "
$ws.Range("D19").Value = "'244.90"
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").Value = "'0.000007763"
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "2.152.79"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'6.902"
$ws.Range("E24").Value = "  -0.73%  "
$ws.Range("D25").Value = "'167.52"
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("D26").Value = "'9.192"
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("E27").Value = "  +10.16%  "
$ws.Range("D28").Value = "'18.84"
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("D29").Value = "'2.035"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").Value = "'1.356"
$ws.Range("E30").Value = "  -2.62%  "
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").Value = "'4.297"
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("D33").Value = "'0.05577"
$ws.Range("E33").Value = "  +4.37%  "
$ws.Range("D34").Value = "'4.061"
$ws.Range("E34").Value = "  +0.86%  "
$ws.Range("D35").Value = "'1.264"
$ws.Range("E35").Value = "  +1.18%  "
$ws.Range("D36").Value = "'0.7290"
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "'0.01931"
$ws.Range("E38").Value = "  +1.04%  "
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").Value = "'0.4404"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").Value = "'5.997"
$ws.Range("E41").Value = "  -2.71%  "
$ws.Range("D42").Value = "'72.04"
$ws.Range("E42").Value = "  -0.25%  "
$ws.Range("D43").Value = "'1.000"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "'0.8355"
$ws.Range("E44").Value = "  +1.35%  "
$ws.Range("D45").Value = "'1.865"
$ws.Range("E45").Value = "  -1.76%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "'7.567"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "'100.46"
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("D48").Value = "'9.716"
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("D49").Value = "'977.72"
$ws.Range("E49").Value = "  +7.93%  "
$ws.Range("D50").Value = "2.057.37"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").Value = "'36.20"
$ws.Range("E51").Value = "  -0.08%  "
